$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J2: report-type code text "001" -> "002" (keep as text, strip the
#     quote-prefix style Excel applies when a leading-zero string is typed) ---
$ws.Range("J2").Value = "'002"
$ws.Range("J2").Style = "Normal"

# --- N2: report date text ---
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# --- Updated numeric figures ---
$ws.Range("O2").Value = 3138321859.39
$ws.Range("P2").Value = 1249765848.17
$ws.Range("Q2").Value = 244120322.42
$ws.Range("S2").Value = 465713337.53
$ws.Range("U2").Value = 664743052.6799999
$ws.Range("W2").Value = 1881273163.17
$ws.Range("X2").Value = 526469615.41
$ws.Range("AB2").Value = 1257048696.22
$ws.Range("AF2").Value = 88.5481090016
$ws.Range("AG2").Value = 59.9451951539

# --- Cells cleared to blank (text) in this update cycle ---
$blankCells = "R2", "T2", "V2", "Y2", "Z2", "AA2", "AC2", "AD2", "AE2"
foreach ($addr in $blankCells) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}
